$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.301.98"
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").Value = "1.585.41"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'209.55"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  -1.01%  "
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("D10").Value = "'19.59"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("D12").Value = "1.807.71"
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("D13").Value = "1.586.89"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").Value = "'4.02"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").Value = "'0.518"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.25%  "
$ws.Range("D16").Value = "'64.45"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.22%  "
$ws.Range("D17").Value = "26.294.21"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").Value = "'7.23"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").Value = "'207.23"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.70%  "
$ws.Range("D22").Value = "'4.26"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("E23").Value = "  -3.91%  "
$ws.Range("E24").Value = "  -1.73%  "
$ws.Range("D25").Value = "'144.54"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("D27").Value = "'7.01"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.40%  "
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("D29").Value = "'15.27"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.91%  "
$ws.Range("D30").Value = "'0.0504"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.10%  "
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("E32").Value = "  -0.84%  "
$ws.Range("E33").Value = "  +13.22%  "
$ws.Range("E34").Value = "  -1.15%  "
$ws.Range("D35").Value = "1.281.46"
$ws.Range("E35").Value = "  -1.31%  "
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("D37").Value = "'0.614"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("E38").Value = "  -1.18%  "
$ws.Range("E39").Value = "  -1.31%  "
$ws.Range("D40").Value = "'0.819"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.45%  "
$ws.Range("D41").Value = "'5.47"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.99%  "
$ws.Range("E42").Value = "  -1.68%  "
$ws.Range("D43").Value = "'2.14"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.11%  "
$ws.Range("D44").Value = "'62.31"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.48%  "
$ws.Range("D45").Value = "1.720.53"
$ws.Range("E45").Value = "  -0.88%  "
$ws.Range("D46").Value = "'88.87"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("D47").Value = "'1.56"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("E48").Value = "  +0.56%  "
$ws.Range("E49").Value = "  -1.63%  "
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("D51").Value = "'7.46"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.51%  "
